$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.091.28'
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").Value = '3.518.38'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '3.518.50'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.10'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.376'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.19%  '
$ws.Range("D13").Value = '4.120.39'
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.75%  '
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000178'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.517.65'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '64.150.36'
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.99%  '
$ws.Range("E21").Value = '  -0.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '382.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.09%  '
$ws.Range("D23").Value = '3.662.21'
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.569'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  +3.33%  '
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.84%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.37'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("D34").Value = '3.535.17'
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  -2.22%  '
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.41'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '159.37'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0787'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.811'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.03%  '
$ws.Range("E47").Value = '  -2.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("E49").Value = '  -2.66%  '
$ws.Range("D50").Value = '2.480.77'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("E51").Value = '  -0.90%  '
